# The "Partner Match" sheet had a few accidental duplicate partner rows
# (consecutive rows in column A repeating the same partner name). Remove
# those duplicate rows so the list is de-duplicated; rows below each
# deleted row shift up automatically, so the used range shrinks from
# A1:D117 to A1:D113.
#
# Duplicate rows (original row numbers, before any deletion):
#   20  Capital One                (dup of row 19)
#   42  Firefly Innovations        (dup of row 41)
#   76  Nonfiction Design          (dup of row 75)
#   95  Strada Education Network   (dup of row 94)
#
# Delete from the bottom up so earlier row numbers stay valid while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$duplicateRows = @(95, 76, 42, 20)
foreach ($r in $duplicateRows) {
    $ws.Rows.Item($r).Delete()
}
